$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3561.3845
$ws.Range("I51").Value = 2997.5
$ws.Range("J51").Value = 3812
$ws.Range("K51").Value = 2997.5
$ws.Range("L51").Value = 3812
$ws.Range("M51").Value = -2513.5
$ws.Range("N51").Value = -4780

$ws.Range("H86").Value = 2653.2727
$ws.Range("I86").Value = 2648.875
$ws.Range("J86").Value = 2665
$ws.Range("K86").Value = 2648.875
$ws.Range("L86").Value = 2665
$ws.Range("M86").Value = -1525.875
$ws.Range("N86").Value = -4911

$ws.Range("H89").Value = 2653.2727
$ws.Range("I89").Value = 2648.875
$ws.Range("J89").Value = 2665
$ws.Range("K89").Value = 13244.375
$ws.Range("L89").Value = 13325
$ws.Range("M89").Value = -7628.375
$ws.Range("N89").Value = -24557

$ws.Range("H113").Value = 2896.625
$ws.Range("I113").Value = 2217.6667
$ws.Range("J113").Value = 3304
$ws.Range("K113").Value = 2217.6667
$ws.Range("L113").Value = 3304
$ws.Range("M113").Value = 1036.3333
$ws.Range("N113").Value = -9812

$ws.Range("H121").Value = 1954.375
$ws.Range("J121").Value = 2306.6667
$ws.Range("L121").Value = 6920.000100000001
$ws.Range("N121").Value = -10414.0001

$ws.Range("H129").Value = 1128.7046
$ws.Range("J129").Value = 1212.1794
$ws.Range("L129").Value = 3636.5382
$ws.Range("N129").Value = -13636.5382

$ws.Range("H132").Value = 2259.3555
$ws.Range("I132").Value = 2250.7778
$ws.Range("J132").Value = 2272.2222
$ws.Range("K132").Value = 6752.3334
$ws.Range("L132").Value = 6816.6666
$ws.Range("M132").Value = -4222.3334
$ws.Range("N132").Value = -11876.6666

$ws.Range("H137").Value = 1303.4595
$ws.Range("I137").Value = 1230.5938
$ws.Range("J137").Value = 1769.8
$ws.Range("K137").Value = 3691.7814
$ws.Range("L137").Value = 5309.4
$ws.Range("M137").Value = -1141.7814
$ws.Range("N137").Value = -10409.4

$ws.Range("H138").Value = 2249.4736
$ws.Range("I138").Value = 1298.0513
$ws.Range("J138").Value = 3252.3242
$ws.Range("K138").Value = 3894.1539
$ws.Range("L138").Value = 9756.972600000001
$ws.Range("M138").Value = 1245.8461
$ws.Range("N138").Value = -20036.9726

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1583.5834
$ws.Range("I2").Value = 1149.25
$ws.Range("J2").Value = 2452.25
$ws.Range("K2").Value = 1149.25
$ws.Range("L2").Value = 2452.25
$ws.Range("M2").Value = -1036.25
$ws.Range("N2").Value = -2678.25

$ws.Range("H32").Value = 12580.085
$ws.Range("I32").Value = 14037.4
$ws.Range("J32").Value = 4483.8887
$ws.Range("K32").Value = 14037.4
$ws.Range("L32").Value = 4483.8887
$ws.Range("M32").Value = -13750.4
$ws.Range("N32").Value = -5057.8887

$ws.Range("H35").Value = 13634.8
$ws.Range("I35").Value = 14043.5
$ws.Range("K35").Value = 14043.5
$ws.Range("M35").Value = -13637.5

$ws.Range("H61").Value = 2543.8
$ws.Range("I61").Value = 2580.1667
$ws.Range("J61").Value = 2489.25
$ws.Range("K61").Value = 2580.1667
$ws.Range("L61").Value = 2489.25
$ws.Range("M61").Value = -2368.1667
$ws.Range("N61").Value = -2913.25

$ws.Range("H109").Value = 38249.875
$ws.Range("J109").Value = 38249.875
$ws.Range("L109").Value = 38249.875
$ws.Range("N109").Value = -41023.875

$ws.Range("H116").Value = 1583.5834
$ws.Range("I116").Value = 1149.25
$ws.Range("J116").Value = 2452.25
$ws.Range("K116").Value = 1149.25
$ws.Range("L116").Value = 2452.25
$ws.Range("M116").Value = 1144.75
$ws.Range("N116").Value = -7040.25

$ws.Range("H123").Value = 24166.2
$ws.Range("J123").Value = 24166.2
$ws.Range("L123").Value = 24166.2
$ws.Range("N123").Value = -33966.2

$ws.Range("H132").Value = 4204.7236
$ws.Range("I132").Value = 5003.2
$ws.Range("J132").Value = 2795.647
$ws.Range("K132").Value = 15009.6
$ws.Range("L132").Value = 8386.940999999999
$ws.Range("M132").Value = -12479.6
$ws.Range("N132").Value = -13446.941

$ws.Range("H136").Value = 2543.8
$ws.Range("I136").Value = 2580.1667
$ws.Range("J136").Value = 2489.25
$ws.Range("K136").Value = 7740.500100000001
$ws.Range("L136").Value = 7467.75
$ws.Range("M136").Value = -5190.500100000001
$ws.Range("N136").Value = -12567.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1583.5834
$ws.Range("I3").Value = 1149.25
$ws.Range("J3").Value = 2452.25
$ws.Range("K3").Value = 1149.25
$ws.Range("L3").Value = 2452.25
$ws.Range("M3").Value = -1035.25
$ws.Range("N3").Value = -2680.25

$ws.Range("H94").Value = 1377.2632
$ws.Range("I94").Value = 671.63635
$ws.Range("J94").Value = 2347.5
$ws.Range("K94").Value = 671.63635
$ws.Range("L94").Value = 2347.5
$ws.Range("M94").Value = -220.63635
$ws.Range("N94").Value = -3249.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9662.5
$ws.Range("I4").Value = 8650
$ws.Range("J4").Value = 10000
$ws.Range("K4").Value = 8650
$ws.Range("M4").Value = -8538
$ws.Range("N4").Value = -10224

$ws.Range("H16").Value = 1252.3077
$ws.Range("I16").Value = 1329.8334
$ws.Range("J16").Value = 1185.8572
$ws.Range("K16").Value = 1329.8334
$ws.Range("L16").Value = 1185.8572
$ws.Range("M16").Value = -1042.8334
$ws.Range("N16").Value = -1759.8572

$ws.Range("H36").Value = 17548
$ws.Range("I36").Value = 17548
$ws.Range("K36").Value = 17548
$ws.Range("M36").Value = -17160

$ws.Range("H40").Value = 17548
$ws.Range("I40").Value = 17548
$ws.Range("K40").Value = 17548
$ws.Range("M40").Value = -17388

$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498

$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488

$ws.Range("H94").Value = 1736.8235
$ws.Range("I94").Value = 1100
$ws.Range("J94").Value = 2002.1666
$ws.Range("K94").Value = 1100
$ws.Range("L94").Value = 2002.1666
$ws.Range("M94").Value = -649
$ws.Range("N94").Value = -2904.1666

$ws.Range("H97").Value = 38399.4
$ws.Range("J97").Value = 38399.4
$ws.Range("L97").Value = 38399.4
$ws.Range("N97").Value = -40381.4

$ws.Range("H113").Value = 1252.3077
$ws.Range("I113").Value = 1329.8334
$ws.Range("J113").Value = 1185.8572
$ws.Range("K113").Value = 1329.8334
$ws.Range("L113").Value = 1185.8572
$ws.Range("M113").Value = 840.1666
$ws.Range("N113").Value = -5525.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1489.625
$ws.Range("I132").Value = 1076.5
$ws.Range("J132").Value = 1902.75
$ws.Range("K132").Value = 9688.5
$ws.Range("L132").Value = 17124.75
$ws.Range("M132").Value = -7158.5
$ws.Range("N132").Value = -22184.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 9297.294
$ws.Range("J109").Value = 9297.294
$ws.Range("L109").Value = 9297.294
$ws.Range("N109").Value = -11377.294

$ws.Range("H122").Value = 2856
$ws.Range("I122").Value = 2181.4
$ws.Range("J122").Value = 3337.8572
$ws.Range("K122").Value = 6544.200000000001
$ws.Range("L122").Value = 10013.5716
$ws.Range("M122").Value = -4094.200000000001
$ws.Range("N122").Value = -14913.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H121").Value = 45059.84
$ws.Range("J121").Value = 45059.84
$ws.Range("L121").Value = 45059.84
$ws.Range("N121").Value = -48553.84

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 50000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = $null
